$excel.Calculation = -4135
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Drive Team Data")

$ws.Range("K5").Value2 = 59.09564474825818
$ws.Range("L5").Value2 = 27.370449678300258
$ws.Range("M5").Value2 = 16.565002471134697
$ws.Range("N5").Value2 = 3.0637701471444605
$ws.Range("O5").Value2 = 4.917836713430747
$ws.Range("P5").Value2 = 0.428170988057082
$ws.Range("Q5").Value2 = 0.0036401456082147607
$ws.Range("R5").Value2 = 12.387543253966445
$ws.Range("S5").Value2 = -1.3321947054345664
$ws.Range("T5").Value2 = 101.37244838466259
$ws.Range("K6").Value2 = 101.450000000221
$ws.Range("K7").Value2 = 71.63157894766981
$ws.Range("T7").Value2 = 70.88354235810294
$ws.Range("K8").Value2 = 71.35483870988972
$ws.Range("L8").Value2 = 32.80645161230885
$ws.Range("M8").Value2 = 10.304347825212854
$ws.Range("O8").Value2 = 3.4999999999999996
$ws.Range("S8").Value2 = -2.096774193406868
$ws.Range("T8").Value2 = 97.1023898048047
$ws.Range("K9").Value2 = 70.6513470692239
$ws.Range("L9").Value2 = 45.055599682434995
$ws.Range("M9").Value2 = 15.833200953251318
$ws.Range("N9").Value2 = 2.4244004171083042
$ws.Range("O9").Value2 = 4.642335766448563
$ws.Range("P9").Value2 = 0.11366006258184523
$ws.Range("Q9").Value2 = 0.3660062565111129
$ws.Range("R9").Value2 = 10.285941223928274
$ws.Range("S9").Value2 = -1.9651347067267961
$ws.Range("T9").Value2 = 114.50229143613548
$ws.Range("K11").Value2 = 70.3877068553645
$ws.Range("K13").Value2 = 63.96424581064444
$ws.Range("L13").Value2 = 27.432432432728802
$ws.Range("M13").Value2 = 13.049549549105542
$ws.Range("N13").Value2 = 2.054054054098758
$ws.Range("O13").Value2 = 4.103603603510607
$ws.Range("P13").Value2 = 0.22522522521749858
$ws.Range("Q13").Value2 = 0.11261261260874929
$ws.Range("R13").Value2 = 23.040540540987582
$ws.Range("T13").Value2 = 81.45043734013018
$ws.Range("K16").Value2 = 66.27368421051614
$ws.Range("L16").Value2 = 26.68421052541352
$ws.Range("M16").Value2 = 11.946428571336416
$ws.Range("N16").Value2 = 1.9310344831925486
$ws.Range("O16").Value2 = 3.6551724143213637
$ws.Range("P16").Value2 = 0.3103448275484741
$ws.Range("R16").Value2 = 20.625000003604914
$ws.Range("S16").Value2 = -4.105263157563213
$ws.Range("T16").Value2 = 122.00496218651159
$ws.Range("K24").Value2 = 63.201149425670074
$ws.Range("L24").Value2 = 29.02298850589806
$ws.Range("M24").Value2 = 12.741379310522712
$ws.Range("O24").Value2 = 4.2873563218498605
$ws.Range("P24").Value2 = 0.14367816092493066
$ws.Range("Q24").Value2 = 0.14367816092493066
$ws.Range("R24").Value2 = 21.436781609249305
$ws.Range("T24").Value2 = 76.89834309861402
$ws.Range("K31").Value2 = 70.60762510032772
$ws.Range("L31").Value2 = 45.055599682434995
$ws.Range("M31").Value2 = 15.833200953251318
$ws.Range("N31").Value2 = 2.4244004171083042
$ws.Range("O31").Value2 = 4.642335766448563
$ws.Range("P31").Value2 = 0.11366006258184523
$ws.Range("Q31").Value2 = 0.3660062565111129
$ws.Range("R31").Value2 = 10.285941223928274
$ws.Range("S31").Value2 = -1.9698173152457137
$ws.Range("T31").Value2 = 114.55540918708306
$ws.Range("L33").Value2 = 29.999999999999996
$ws.Range("K35").Value2 = 53.39285713999491
$ws.Range("L35").Value2 = 9.285714288004081
$ws.Range("M35").Value2 = 19.642857138659185
$ws.Range("O35").Value2 = 2.6785714279989796
$ws.Range("P35").Value2 = 0.8928571426663265
$ws.Range("R35").Value2 = 24.46428571333163
$ws.Range("T35").Value2 = 119.65465824202485
$ws.Range("K38").Value2 = 55.647281921267705
$ws.Range("L38").Value2 = 27.701587301341245
$ws.Range("M38").Value2 = 17.878730158582986
$ws.Range("N38").Value2 = 3.429101019433519
$ws.Range("O38").Value2 = 5.196825396768324
$ws.Range("P38").Value2 = 0.4661723818217156
$ws.Range("Q38").Value2 = 0.004444444447668148
$ws.Range("R38").Value2 = 10.04444444447668
$ws.Range("T38").Value2 = 89.23197499716709
$ws.Range("K39").Value2 = 101.450000000221
$ws.Range("K41").Value2 = 71.35483870988972
$ws.Range("L41").Value2 = 32.80645161230885
$ws.Range("M41").Value2 = 10.304347825212854
$ws.Range("O41").Value2 = 3.4999999999999996
$ws.Range("S41").Value2 = -2.096774193406868
$ws.Range("T41").Value2 = 97.1023898048047
$ws.Range("K46").Value2 = 83.67816092039192
$ws.Range("T46").Value2 = 67.46196313619156
